$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Data")

# New rows of arrival data to append to the bottom of the table.
$newRows = @(
    @{ Number = 182; Date = "Sunday, Jan 15"; Time = "4:55 PM"; Flight = "FR4059"; From = "Malaga"; Short = "(AGP)"; Airline = "Buzz "; Model = "B38M"; AircraftId = "(SP-RZE)"; Status = "4:41 PM"; Difference = "0 hours, -14 minutes" },
    @{ Number = 183; Date = "Sunday, Jan 15"; Time = "5:10 PM"; Flight = "FR1021"; From = "London"; Short = "(STN)"; Airline = "Ryanair "; Model = "B738"; AircraftId = "(EI-EBP)"; Status = "5:28 PM"; Difference = "0 hours, 18 minutes" },
    @{ Number = 184; Date = "Sunday, Jan 15"; Time = "5:13 PM"; Flight = "UNKNOWN"; From = "London"; Short = "(LTN)"; Airline = "NetJets Europe "; Model = "F2TH"; AircraftId = "(CS-DLF)"; Status = "5:09 PM"; Difference = "0 hours, -4 minutes" }
)

$startRow = 183
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.Number
    $ws.Cells.Item($r, 2).Value = $row.Date
    $ws.Cells.Item($r, 3).Value = $row.Time
    $ws.Cells.Item($r, 4).Value = $row.Flight
    $ws.Cells.Item($r, 5).Value = $row.From
    $ws.Cells.Item($r, 6).Value = $row.Short
    $ws.Cells.Item($r, 7).Value = $row.Airline
    $ws.Cells.Item($r, 8).Value = $row.Model
    $ws.Cells.Item($r, 9).Value = $row.AircraftId
    $ws.Cells.Item($r, 10).Value = $row.Status
    $ws.Cells.Item($r, 11).ClearFormats()
    $ws.Cells.Item($r, 12).Value = $row.Difference
    $ws.Cells.Item($r, 13).ClearFormats()
}
